$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = 249335
$ws.Range("E3").Value = 1036488011

# Row 6
$ws.Range("C6").Value = 20796
$ws.Range("E6").Value = 360718555

# Row 7
$ws.Range("C7").Value = 7018
$ws.Range("E7").Value = 291118761

# Row 53
$ws.Range("C53").Value = 141687
$ws.Range("E53").Value = 590077232

# Row 82
$ws.Range("C82").Value = 8453
$ws.Range("E82").Value = 124867404

# Row 92
$ws.Range("C92").Value = 409316
$ws.Range("E92").Value = 1597549233

# Row 93
$ws.Range("C93").Value = 209672
$ws.Range("E93").Value = 1310272697

# Row 94
$ws.Range("C94").Value = 94241
$ws.Range("E94").Value = 919287786

# Row 95
$ws.Range("C95").Value = 50809
$ws.Range("E95").Value = 934729383

# Row 96
$ws.Range("C96").Value = 17327
$ws.Range("E96").Value = 797769865

# Row 166
$ws.Range("C166").Value = 35931
$ws.Range("E166").Value = 210610514
